$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("M2").Value = 8.949653
$ws.Range("N2").Value = 17.899306
$ws.Range("O2").Value = 0.1668927877080592
$ws.Range("P2").Value = 0.1610232428880788
$ws.Range("Q2").Value = 1.550387171019667
$ws.Range("R2").Value = 9.302323026118
$ws.Range("S2").Value = 0.1668927877080592
$ws.Range("T2").Value = 0.1610232428880788
$ws.Range("O3").Value = 0.06249788578732534
$ws.Range("P3").Value = 0.09044979457765322
$ws.Range("S3").Value = 0.06249788578732534
$ws.Range("T3").Value = 0.09044979457765322
$ws.Range("M4").Value = 0.050758
$ws.Range("N4").Value = 0.152274
$ws.Range("O4").Value = 0.0009465332475444208
$ws.Range("P4").Value = 0.001369866143834812
$ws.Range("Q4").Value = 0.008793028291333334
$ws.Range("R4").Value = 0.079137254622
$ws.Range("S4").Value = 0.0009465332475444208
$ws.Range("T4").Value = 0.001369866143834812
$ws.Range("M5").Value = 40.766071
$ws.Range("N5").Value = 81.532142
$ws.Range("O5").Value = 0.7602041367519689
$ws.Range("P5").Value = 0.7334680967212543
$ws.Range("Q5").Value = 7.062083132304333
$ws.Range("R5").Value = 42.372498793826
$ws.Range("S5").Value = 0.7602041367519689
$ws.Range("T5").Value = 0.7334680967212543
$ws.Range("M6").Value = 0.3317233333333334
$ws.Range("N6").Value = 0.9951700000000001
$ws.Range("O6").Value = 0.006185964064507279
$ws.Range("P6").Value = 0.008952609705925435
$ws.Range("Q6").Value = 0.05746587050111112
$ws.Range("R6").Value = 0.5171928345100001
$ws.Range("S6").Value = 0.006185964064507279
$ws.Range("T6").Value = 0.008952609705925435
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1754986666666667
$ws.Range("N7").Value = 0.526496
$ws.Range("O7").Value = 0.003272692440594897
$ws.Range("P7").Value = 0.004736389963253432
$ws.Range("Q7").Value = 0.03040239452088889
$ws.Range("R7").Value = 0.273621550688
$ws.Range("S7").Value = 0.003272692440594897
$ws.Range("T7").Value = 0.004736389963253432
